$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header column Z1 = "OMUSDT" ---------------------------------
# Copy the header formatting (bold font + border) from Y1 so the new
# shared-string entry gets the same style index (s="1") instead of an
# unformatted cell.
$ws.Range("Y1").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Range("Z1").Value = "OMUSDT"

# --- New data rows 55:57 ------------------------------------------------
# Column A carries the date style (s="2") on every existing data row;
# copy that formatting down onto the three new rows before filling values.
$ws.Range("A54").Copy()
$ws.Range("A55:A57").PasteSpecial(-4122)

# Row 55 (2024-07-12)
$ws.Cells.Item(55,1).Value = [double]"45485"
$ws.Cells.Item(55,2).Value = [double]"617.864257684"
$ws.Cells.Item(55,3).Value = [double]"217.291800538"
$ws.Cells.Item(55,4).Value = [double]"0"
$ws.Cells.Item(55,5).Value = [double]"13.62810763023"
$ws.Cells.Item(55,6).Value = [double]"0"
$ws.Cells.Item(55,7).Value = [double]"113.21935185"
$ws.Cells.Item(55,8).Value = [double]"0"
$ws.Cells.Item(55,9).Value = [double]"204.490246811"
$ws.Cells.Item(55,10).Value = [double]"0"
$ws.Cells.Item(55,11).Value = [double]"21.4493759491004"
$ws.Cells.Item(55,12).Value = [double]"0"
$ws.Cells.Item(55,13).Value = [double]"0"
$ws.Cells.Item(55,14).Value = [double]"137.6900122272"
$ws.Cells.Item(55,15).Value = [double]"54.220932506"
$ws.Cells.Item(55,16).Value = [double]"0"
$ws.Cells.Item(55,17).Value = [double]"2.0856E-06"
$ws.Cells.Item(55,18).Value = [double]"0"
$ws.Cells.Item(55,19).Value = [double]"0"
$ws.Cells.Item(55,20).Value = [double]"0"
$ws.Cells.Item(55,21).Value = [double]"316.2405328225512"
$ws.Cells.Item(55,22).Value = [double]"0"
$ws.Cells.Item(55,23).Value = [double]"0"
$ws.Cells.Item(55,24).Value = [double]"0"
$ws.Cells.Item(55,25).Value = [double]"0"

# Row 56 (2024-07-13)
$ws.Cells.Item(56,1).Value = [double]"45486"
$ws.Cells.Item(56,2).Value = [double]"631.8987144248"
$ws.Cells.Item(56,3).Value = [double]"220.2073940555"
$ws.Cells.Item(56,4).Value = [double]"0"
$ws.Cells.Item(56,5).Value = [double]"13.35302387655"
$ws.Cells.Item(56,6).Value = [double]"0"
$ws.Cells.Item(56,7).Value = [double]"121.8939588"
$ws.Cells.Item(56,8).Value = [double]"0"
$ws.Cells.Item(56,9).Value = [double]"208.073670123"
$ws.Cells.Item(56,10).Value = [double]"0"
$ws.Cells.Item(56,11).Value = [double]"21.7927250151774"
$ws.Cells.Item(56,12).Value = [double]"0"
$ws.Cells.Item(56,13).Value = [double]"0"
$ws.Cells.Item(56,14).Value = [double]"159.76731589088"
$ws.Cells.Item(56,15).Value = [double]"54.078726988"
$ws.Cells.Item(56,16).Value = [double]"0"
$ws.Cells.Item(56,17).Value = [double]"2.0832E-06"
$ws.Cells.Item(56,18).Value = [double]"0"
$ws.Cells.Item(56,19).Value = [double]"0"
$ws.Cells.Item(56,20).Value = [double]"0"
$ws.Cells.Item(56,21).Value = [double]"305.3665662813227"
$ws.Cells.Item(56,22).Value = [double]"0"
$ws.Cells.Item(56,23).Value = [double]"0"
$ws.Cells.Item(56,24).Value = [double]"0"
$ws.Cells.Item(56,25).Value = [double]"0"

# Row 57 (2024-07-14) -- sparse row, several columns legitimately stay blank
$ws.Cells.Item(57,1).Value = [double]"45487"
$ws.Cells.Item(57,3).Value = [double]"225.002002658"
$ws.Cells.Item(57,4).Value = [double]"0"
$ws.Cells.Item(57,5).Value = [double]"14.29289336829"
$ws.Cells.Item(57,7).Value = [double]"125.08513545"
$ws.Cells.Item(57,8).Value = [double]"0"
$ws.Cells.Item(57,10).Value = [double]"0"
$ws.Cells.Item(57,11).Value = [double]"21.5047660415346"
$ws.Cells.Item(57,12).Value = [double]"0"
$ws.Cells.Item(57,13).Value = [double]"0"
$ws.Cells.Item(57,14).Value = [double]"159.13079272"
$ws.Cells.Item(57,16).Value = [double]"0"
$ws.Cells.Item(57,17).Value = [double]"2.208E-06"
$ws.Cells.Item(57,20).Value = [double]"0"
$ws.Cells.Item(57,24).Value = [double]"0"
$ws.Cells.Item(57,26).Value = [double]"209.90145507716"
